$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "2025-08-06 01:22:53"
$ws.Range("B13").Value = "remove-repo"
$ws.Range("C13").Value = "new-organization97"
$ws.Range("D13").Value = "Devops"
$ws.Range("E13").Value = "deerepo"

# "False" would otherwise be auto-coerced to a Boolean by Excel's input
# parser; a leading apostrophe forces literal text (matching the other
# I-column cells in this sheet, which are plain text "False"/"True"
# values, not booleans). The quote-prefix nudges the cell's style, so
# reset it back to Normal afterwards to leave the cell unstyled like its
# neighbours.
$ws.Range("I13").Value = "'False"
$ws.Range("I13").Style = "Normal"
